$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 403, shifting existing rows 403:497 down to 404:498
$ws.Rows("403:403").Insert()

# Populate the newly inserted row 403 with the new data record
$ws.Range("A403").Value = 4
$ws.Range("B403").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C403").Value = "Los Lagos"
$ws.Range("D403").Value = 44782
$ws.Range("E403").Value = 10
$ws.Range("F403").Value = 100112006
$ws.Range("G403").Value = "Repollo"
$ws.Range("H403").Value = "Crespo record"
$ws.Range("I403").Value = "Primera"
$ws.Range("J403").Value = 1400
$ws.Range("K403").Value = 2000
$ws.Range("L403").Value = 2000
$ws.Range("M403").Value = 2000
$ws.Range("N403").Value = "$/unidad"
$ws.Range("O403").Value = "Región Metropolitana"
$ws.Range("P403").Value = 2000
$ws.Range("Q403").Value = 1
$ws.Range("R403").Value = "Hortaliza"
